# Daily attendance processing - 2025-10-14 23:20:41
# Swap the order of the two comma-separated "Recorded By" entries
# (column G) for the specific rows touched by this run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(3,6,7,10,11,12,13,14,15,17,18,19,30,33,34,37,38,39,40,41,42,44,45,46,57,60,61,64,65,66,67,68,69,71,72,73,86,87,88,89,90,93,95,96,97,112,113,114,115,116,119,121,122,123,138,139,140,141,142,145,147,148,149)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $current = [string]$cell.Value2
    $parts = $current -split ", ", 2
    if ($parts.Length -eq 2) {
        $swapped = $parts[1] + ", " + $parts[0]
        $cell.Value = $swapped
    }
}
